$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Piroxicam")

# Update the RxCUI code for Piroxicam in D2 from text "C0031990" to the numeric value 8356
$ws.Range("D2").Value = 8356

# Update the active selection to match the saved view state
$ws.Range("C10").Select()
